# Actualización automática 2025-06-25 17:30:09
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Update VENTA (D2) and the dependent totals/percentages on the sheet.
$ws.Range("D2").Value = 4098.58
$ws.Range("E2").Value = -4098.58

$ws.Range("D4").Value = 4647.04
$ws.Range("E4").Value = 12852.96
$ws.Range("F4").Value = 0.2655451428571429
